$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1; this pushes the existing header row
# (currently row 1: A/B/C/D with the bold+centered style) down to row 2,
# and every data row below it shifts down by one as well.
$ws.Rows.Item(1).Insert()

# Row 2 (the old header row) gets " Annual" appended to the year labels,
# and the account-name label moves from the old A1 into the new A2.
$ws.Range("B2").Value = "2024 Annual"
$ws.Range("C2").Value = "2023 Annual"
$ws.Range("D2").Value = "2022 Annual"

# New row 1: a single note spanning columns A:D, styled italic + red.
$ws.Range("A1").Value = "Note: The date header (Row 2) supports: '2023 Annual', '2023 Q1', '2023-01'"
$ws.Range("A1:D1").Merge()
$ws.Range("A1").Font.Italic = $true
$ws.Range("A1").Font.Color = 255
